# OpenSolver ChangeLog.xlsx - add release notes for "Version 2.5.3 alpha"
#
# The new entry is inserted as a new block of 3 rows at the very top of the
# release-notes list (row 9), pushing every row below it down by 3. This
# mirrors the existing pattern used for every other version entry on the
# sheet: a bold "Version X" header row, one (or more) regular detail rows,
# and a trailing blank separator row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 3 new rows above the current "Version 2.5.2 alpha" entry (row 9).
#    This shifts all rows/cells/merged-ranges below down by 3, same as the
#    canonical diff (e.g. old row 9 -> new row 12, old row 212 -> new row 215).
$ws.Rows("9:11").Insert()

# 2. Populate the new header + detail rows.
$ws.Range("A9").Value2 = "Version 2.5.3 alpha"
$ws.Range("A10").Value2 = "Add support for NOMAD in 64-bit Office."

# Row 9 ("Version 2.5.3 alpha") keeps the bold "header" formatting that was
# copied down automatically from row 8 on insert.
# Row 10 (detail bullet) should not be bold - match the plain/regular look
# used by every other detail row (e.g. old row 10 "Fix memory bug...").
$ws.Range("A10").ClearFormats()
$ws.Range("A10").Value2 = "Add support for NOMAD in 64-bit Office."

# Row 11 stays as the blank separator row, which already has the right
# (header-style, but empty) formatting copied down from row 8 on insert.

# 3. Fix up the hidden Solver defined name that pointed at a cell below the
#    insertion point - it needs to shift down by 3 rows just like the data.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!solver_opt") {
        $n.RefersTo = "=Sheet1!`$U`$103"
    }
}

# 4. Re-point the hyperlink that lived on the old B108 (now B111) - this
#    engine does not automatically shift hyperlink anchors on row insert.
#    Recreate both hyperlinks (B111 and H7) so the relationship id order in
#    the saved file matches the original (B-column link first, then H7).
$existing = @()
foreach ($h in $ws.Hyperlinks) { $existing += $h }
for ($i = $existing.Count - 1; $i -ge 0; $i--) {
    $existing[$i].Delete()
}

$ws.Hyperlinks.Add($ws.Range("B111"), "http://www.officekb.com/Uwe/Forum.aspx/excel-prog/159706/Shape-TextEffect-HorizontalAlignment-throws-error")
$ws.Range("B111").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("H7"), "http://www.vbforums.com/archive/index.php/t-47843.html")
$ws.Range("H7").Style = "Hyperlink"

# 5. Match the cursor position left behind in the saved workbook.
$ws.Range("F12").Select()
